# Fixed a bug in symbolModifier
# Rewrites the symbol/reel data rows (rows 2-21) on the active sheet with
# the corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(401,  9, 48, 67, 75, 45),
    @(601,  9, 60, 67, 60, 42),
    @(201,  9, 30, 15, 45, 30),
    @(902,  1,  0,  0,  0,  0),
    @(501,  9, 52, 30, 75, 45),
    @(1202, 2, 10, 10, 10, 10),
    @(101,  9, 30, 15, 60, 15),
    @(901, 16, 15, 45, 60, 60),
    @(801,  3, 67, 65, 52, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(1001,18, 30, 75, 60, 72),
    @(301,  6, 45, 30, 60, 45),
    @(701,  3, 90, 45, 97, 15),
    @(2,    0,  2,  2,  2,  2),
    @(3,    0,  3,  3,  3,  3),
    @(1,    0,  2,  2,  2,  2),
    @(802,  0,  4,  5,  4,  0),
    @(502,  0,  4,  0,  0,  0),
    @(1101, 0, 15, 30, 30,  0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
